$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 691.6667
$ws.Range("I18").Value = 537.5
$ws.Range("K18").Value = 537.5
$ws.Range("M18").Value = -253.5

$ws.Range("H63").Value = 9950
$ws.Range("J63").Value = 9950
$ws.Range("L63").Value = 9950
$ws.Range("N63").Value = -11198

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H66").Value = 9950
$ws.Range("J66").Value = 9950
$ws.Range("L66").Value = 29850
$ws.Range("N66").Value = -36090

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H112").Value = 1686
$ws.Range("I112").Value = 675
$ws.Range("J112").Value = 1878.5714
$ws.Range("K112").Value = 2025
$ws.Range("L112").Value = 5635.7142
$ws.Range("M112").Value = -917
$ws.Range("N112").Value = -7851.7142

$ws.Range("H137").Value = 1223.909
$ws.Range("I137").Value = 1000.1429
$ws.Range("J137").Value = 1615.5
$ws.Range("K137").Value = 3000.4287
$ws.Range("L137").Value = 4846.5
$ws.Range("M137").Value = -450.4287000000004
$ws.Range("N137").Value = -9946.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2321.111
$ws.Range("I2").Value = 1115.4
$ws.Range("K2").Value = 1115.4
$ws.Range("M2").Value = -1002.4

$ws.Range("H32").Value = 335321.97
$ws.Range("I32").Value = 2644.5974
$ws.Range("J32").Value = 3181561.8
$ws.Range("K32").Value = 2644.5974
$ws.Range("L32").Value = 3181561.8
$ws.Range("M32").Value = -2357.5974
$ws.Range("N32").Value = -3182135.8

$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H74").Value = 645.19354
$ws.Range("I74").Value = 579.26086
$ws.Range("J74").Value = 834.75
$ws.Range("K74").Value = 579.26086
$ws.Range("L74").Value = 834.75
$ws.Range("M74").Value = 294.73914
$ws.Range("N74").Value = -2582.75

$ws.Range("H77").Value = 645.19354
$ws.Range("I77").Value = 579.26086
$ws.Range("J77").Value = 834.75
$ws.Range("K77").Value = 2896.3043
$ws.Range("L77").Value = 4173.75
$ws.Range("M77").Value = 1471.6957
$ws.Range("N77").Value = -12909.75

$ws.Range("H116").Value = 2321.111
$ws.Range("I116").Value = 1115.4
$ws.Range("K116").Value = 1115.4
$ws.Range("M116").Value = 1178.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2321.111
$ws.Range("I3").Value = 1115.4
$ws.Range("K3").Value = 1115.4
$ws.Range("M3").Value = -1001.4

$ws.Range("H105").Value = 1713.2963
$ws.Range("I105").Value = 1700.7142
$ws.Range("J105").Value = 1726.8462
$ws.Range("K105").Value = 1700.7142
$ws.Range("L105").Value = 1726.8462
$ws.Range("M105").Value = 46.28580000000011
$ws.Range("N105").Value = -5220.8462

$ws.Range("H134").Value = 5108.82
$ws.Range("I134").Value = 1653.5116
$ws.Range("K134").Value = 4960.5348
$ws.Range("M134").Value = -2425.5348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 535.625
$ws.Range("I19").Value = 212.14285
$ws.Range("J19").Value = 2800
$ws.Range("K19").Value = 212.14285
$ws.Range("L19").Value = 2800
$ws.Range("M19").Value = -42.14285000000001
$ws.Range("N19").Value = -3140

$ws.Range("H22").Value = 596
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -50

$ws.Range("H24").Value = 535.625
$ws.Range("I24").Value = 212.14285
$ws.Range("J24").Value = 2800
$ws.Range("K24").Value = 212.14285
$ws.Range("L24").Value = 2800
$ws.Range("M24").Value = -42.14285000000001
$ws.Range("N24").Value = -3140

$ws.Range("H31").Value = 2060.86
$ws.Range("I31").Value = 1564.075
$ws.Range("J31").Value = 4048
$ws.Range("K31").Value = 1564.075
$ws.Range("L31").Value = 4048
$ws.Range("M31").Value = -1269.075
$ws.Range("N31").Value = -4638

$ws.Range("H34").Value = 2060.86
$ws.Range("I34").Value = 1564.075
$ws.Range("J34").Value = 4048
$ws.Range("K34").Value = 1564.075
$ws.Range("L34").Value = 4048
$ws.Range("M34").Value = -1362.075
$ws.Range("N34").Value = -4452

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 933.4074000000001
$ws.Range("I5").Value = 948.28
$ws.Range("J5").Value = 747.5
$ws.Range("K5").Value = 2844.84
$ws.Range("L5").Value = 2242.5
$ws.Range("M5").Value = -2732.84
$ws.Range("N5").Value = -2466.5

$ws.Range("H34").Value = 1192.1666
$ws.Range("I34").Value = 700
$ws.Range("J34").Value = 1236.909
$ws.Range("K34").Value = 2100
$ws.Range("L34").Value = 3710.727
$ws.Range("M34").Value = -2016
$ws.Range("N34").Value = -3878.727

$ws.Range("H40").Value = 180.94444
$ws.Range("I40").Value = 168.11765
$ws.Range("J40").Value = 399
$ws.Range("K40").Value = 672.4706
$ws.Range("L40").Value = 1596
$ws.Range("M40").Value = -603.4706
$ws.Range("N40").Value = -1734

$ws.Range("H46").Value = 23993.295
$ws.Range("I46").Value = 346.5
$ws.Range("J46").Value = 31269.23
$ws.Range("K46").Value = 1039.5
$ws.Range("L46").Value = 93807.69
$ws.Range("M46").Value = -948.5
$ws.Range("N46").Value = -93989.69

$ws.Range("H61").Value = 289.66666
$ws.Range("J61").Value = 805
$ws.Range("L61").Value = 2415
$ws.Range("N61").Value = -2845

$ws.Range("H107").Value = 430.42856
$ws.Range("I107").Value = 439.2
$ws.Range("J107").Value = 427.6875
$ws.Range("K107").Value = 1317.6
$ws.Range("L107").Value = 1283.0625
$ws.Range("M107").Value = 602.4000000000001
$ws.Range("N107").Value = -5123.0625

$ws.Range("H113").Value = 1136.8667
$ws.Range("I113").Value = 1002.5
$ws.Range("J113").Value = 1146.4642
$ws.Range("K113").Value = 3007.5
$ws.Range("L113").Value = 3439.3926
$ws.Range("M113").Value = -837.5
$ws.Range("N113").Value = -7779.392599999999

$ws.Range("H129").Value = 13335385
$ws.Range("J129").Value = 18520484
$ws.Range("L129").Value = 55561452
$ws.Range("N129").Value = -55571452

$ws.Range("H132").Value = 2091.75
$ws.Range("I132").Value = 1264.091
$ws.Range("J132").Value = 2792.077
$ws.Range("K132").Value = 11376.819
$ws.Range("L132").Value = 25128.693
$ws.Range("M132").Value = -8846.819
$ws.Range("N132").Value = -30188.693

$ws.Range("H135").Value = 933.4074000000001
$ws.Range("I135").Value = 948.28
$ws.Range("J135").Value = 747.5
$ws.Range("K135").Value = 8534.52
$ws.Range("L135").Value = 6727.5
$ws.Range("M135").Value = -5999.52
$ws.Range("N135").Value = -11797.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19404122
$ws.Range("I70").Value = 29611988
$ws.Range("J70").Value = 9178
$ws.Range("K70").Value = 29611988
$ws.Range("L70").Value = 9178
$ws.Range("M70").Value = -29611718
$ws.Range("N70").Value = -9718

$ws.Range("H73").Value = 19404122
$ws.Range("I73").Value = 29611988
$ws.Range("J73").Value = 9178
$ws.Range("K73").Value = 29611988
$ws.Range("L73").Value = 9178
$ws.Range("M73").Value = -29611052
$ws.Range("N73").Value = -11050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1899.7084
$ws.Range("I40").Value = 1881.9412
$ws.Range("K40").Value = 1881.9412
$ws.Range("M40").Value = -1745.9412

$ws.Range("H122").Value = 2066.8
$ws.Range("I122").Value = 2041.091
$ws.Range("J122").Value = 2137.5
$ws.Range("K122").Value = 6123.272999999999
$ws.Range("L122").Value = 6412.5
$ws.Range("M122").Value = -3673.272999999999
$ws.Range("N122").Value = -11312.5
